$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match formatting of existing header cells (bold font, thin border, centered/top aligned)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-23
$iValues = @(8,8,6,6,8,7,9,9,8,6,8,7,7,6,5,6,4,9,5,5,6,4)
$jValues = @(8,8,6,6,8,8,9,9,8,7,8,8,7,7,5,6,4,9,5,5,6,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
